$d = $word.ActiveDocument

function Replace-NextInRange {
    param(
        $Cursor,
        $Limit,
        [string]$Old,
        [string]$New
    )
    $rng = $d.Range($Cursor, $Limit)
    $found = $rng.Find.Execute($Old, $true, $false, $false, $false, $false, $true, 0, $false, $New, 1)
    if (-not $found) {
        Write-Output "NOT FOUND: '$Old' (cursor=$Cursor limit=$Limit)"
        return $Cursor
    }
    return $rng.End
}

# --- Paragraph: "Na figura ?? há um resumo das características dos sítios..." ---
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Na figura ?? há um resumo*") {
        $p1 = $p
        break
    }
}
if ($p1 -ne $null) {
    $limit = $p1.Range.End
    $cursor = $p1.Range.Start
    $cursor = Replace-NextInRange $cursor $limit "??" "2"
}

# --- Paragraph: "Na figura ?? há a taxa U estimada ... fig. ??))... na figura (??)..." ---
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Na figura ?? há a taxa U estimada*") {
        $p2 = $p
        break
    }
}
if ($p2 -ne $null) {
    $limit = $p2.Range.End
    $cursor = $p2.Range.Start
    $cursor = Replace-NextInRange $cursor $limit "??" "3"
    $cursor = Replace-NextInRange $cursor $limit "??" "3"
    $cursor = Replace-NextInRange $cursor $limit "[1]" "May et al. (2012)"
    $cursor = Replace-NextInRange $cursor $limit "??" "4"
    $cursor = Replace-NextInRange $cursor $limit "??" "4"
}

# --- Remove the leftover R console-output paragraph (SourceCode block) ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*danilo*mestrado_Ecologia*") {
        $p.Range.Delete()
        break
    }
}

# --- Rewrite the bibliography entry text/format ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*May F, Giladi I, Ziv Y, Jeltsch F*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $r = $d.Range($target.Range.Start, $target.Range.End)
    $r.Text = "May, F., I. Giladi, Y. Ziv, and F. Jeltsch. 2012. Dispersal and diversity–unifying scale-dependent relationships within the neutral theory. Oikos 121:942–951."
}
